$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - extend existing headers with new column names
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "first_name"
$ws.Range("C1").Value = "last_name"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "car_make"
$ws.Range("F1").Value = "car_model"
$ws.Range("G1").Value = "vin_number"
$ws.Range("H1").Value = "manufactured_date"

# Row 2 - first record
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Allister"
$ws.Range("C2").Value = "Camili"
$ws.Range("D2").Value = "Chegutu"
$ws.Range("E2").Value = "Acura"
$ws.Range("F2").Value = "RL"
$ws.Range("G2").Value = "3434989GDS"
$ws.Range("H2").NumberFormat = "mm-dd-yy"
$ws.Range("H2").Value = "3/3/1990"

# Row 3 - second record
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Kasper"
$ws.Range("C3").Value = "Cogger"
$ws.Range("D3").Value = "Caballococha"
$ws.Range("E3").Value = "Subaru"
$ws.Range("F3").Value = "Leon"
$ws.Range("G3").Value = "WE3942948"
$ws.Range("H3").NumberFormat = "mm-dd-yy"
$ws.Range("H3").Value = "2/28/1984"

# Auto-fit columns to match best-fit widths seen after inserting new data
$ws.Range("A1:H3").Columns.AutoFit()

# Update selection to mirror the authored workbook state
$ws.Range("D11").Select()
